$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Version bump: "0.1" -> "1.0" (cell D2, the "Version: " value).
#    A leading apostrophe forces the numeric-looking text to stay text (not
#    get auto-converted to the number 1); the subsequent PasteSpecial of
#    formats only (copied from the still-original C2, which shares D2's
#    original style) restores the cell's original style index so we don't
#    leave an extra/duplicate style behind just from the text coercion.
$ws.Range("D2").Value = "'1.0"
$ws.Range("C2").Copy() | Out-Null
$ws.Range("D2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# 2) TC1 (row 13) and TC3 (row 32) "Expected Results" cells swap content:
#    TC1's last step now expects the success message, and TC3's last step
#    now expects the failure message (previously it was the other way
#    around).
$ws.Range("D13").Value = "SYSTEM Atualiza os dados bancários do beneficiário na base do RH (SRH); Exibe mensagem de sucesso para o usuário."
$ws.Range("D32").Value = "SYSTEM Identifica que ocorreu uma falha durante a tentativa de atualização dos dados bancários; Mantém os dados consistentes, interrompe a operação; Exibe mensagem de erro (MSG213 - Não foi possível concluir a operação. Falha na comunicação com o sistema de Recursos Humanos) para o usuário."
